$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 281, shifting existing rows 281:305 down to 284:308
$ws.Rows("281:283").Insert()

# New data for the inserted rows (weekly update - newest week prepended)
$newRows = @(
    @{Row=281; D=44783; H="Americana (o)"; I="Primera"; J=100; K=40000; L=42000; M=41000; P=1640},
    @{Row=282; D=44783; H="Americana (o)"; I="Segunda"; J=68;  K=30000; L=32000; M=31118; P=1245},
    @{Row=283; D=44783; H="Inferno";       I="Primera"; J=60;  K=25000; L=26000; M=25500; P=1020}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 2
    $ws.Cells.Item($row, 2).Value  = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = 4
    $ws.Cells.Item($row, 6).Value  = 100112021
    $ws.Cells.Item($row, 7).Value  = "Ají"
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "$/caja 25 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 25
    $ws.Cells.Item($row, 18).Value = "Hortaliza"

    # Column D keeps the date/time custom number format used elsewhere in the sheet
    $ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(280, 4).NumberFormat
}
